# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q2" and "总计", populated
#    with the quarter's fund-holding detail rows.
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (pushing the existing 2021-Q2 summary row down) with the 2022-Q1 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: add the "2022-Q1" sheet right after "2021-Q2"
# ---------------------------------------------------------------------------
$q2sheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q2sheet)
$newSheet.Name = "2022-Q1"

# Re-fetch the "总计" sheet fresh (it shifted from index 2 -> 3 after the insert)
$totalsSheet = $wb.Worksheets.Item(3)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Match the header formatting used on the "总计" sheet's header row (bold,
# centered, top-aligned, thin border).
$totalsSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Data rows (A column index, B fund code, C fund name, D fund scale,
# E stock position total, F position ratio, G held value (100M yuan),
# H position rank). B-G are stored as text, A and H as numbers - matching
# the source data's typed layout.
$data = @(
    @(0, "010088", "工银瑞信优质成长混合A", "19.41", "82.44", "2.24", "0.4348", 9),
    @(1, "012010", "富国泰享回报6个月持有期混合型证券投资基金A", "9.29", "29.91", "1.41", "0.1310", 3),
    @(2, "010089", "工银瑞信优质成长混合C", "1.34", "82.44", "2.24", "0.0300", 9),
    @(3, "161124", "易方达香港恒生综合小型股指数（QDII-LOF）A", "0.28", "92.62", "1.55", "0.0043", 5),
    @(4, "012011", "富国泰享回报6个月持有期混合型证券投资基金C", "0.09", "29.91", "1.41", "0.0013", 3),
    @(5, "006263", "易方达香港恒生综合小型股指数（QDII-LOF）C", "0.06", "92.62", "1.55", "0.0009", 5)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Range("A" + $r).Value = $row[0]
    $newSheet.Range("B" + $r).Value = "'" + $row[1]
    $newSheet.Range("C" + $r).Value = "'" + $row[2]
    $newSheet.Range("D" + $r).Value = "'" + $row[3]
    $newSheet.Range("E" + $r).Value = "'" + $row[4]
    $newSheet.Range("F" + $r).Value = "'" + $row[5]
    $newSheet.Range("G" + $r).Value = "'" + $row[6]
    $newSheet.Range("H" + $r).Value = $row[7]
    $r = $r + 1
}

# Column A on the data rows uses the same bold/border style as the "总计"
# sheet's "A" column (copy it down from the header cell which already has it).
$newSheet.Range("B1").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of "总计"'s data
# ---------------------------------------------------------------------------
$totalsSheet.Rows("2:2").Insert()

$totalsSheet.Range("A2").Value = 0
$totalsSheet.Range("B2").Value = "2022-Q1"
$totalsSheet.Range("C2").Value = 6
$totalsSheet.Range("D2").Value = 0.6

# The pushed-down 2021-Q2 row is now row 3; its running index needs to
# advance from 0 to 1 to stay in sequence with the new row above it.
$totalsSheet.Range("A3").Value = 1

# The row-insert leaves a blank inherited style on B2:D2 and no style on A2;
# clear it and re-copy the real "总计" row style (from the row pushed to 3).
$totalsSheet.Range("A2:D2").ClearFormats()
$totalsSheet.Range("A3").Copy()
$totalsSheet.Range("A2").PasteSpecial(-4122)

# Keep "2021-Q2" the active/selected sheet, as it was before these edits.
$wb.Worksheets.Item(1).Activate()
